$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose NEW value is a plain number-looking string ("268.80", "9.70",
# "19.00", ...). The source workbook stores every Price cell as literal text
# (t="inlineStr"), so a bare $range.Value = "9.70" would let Excel coerce it
# to the number 9.7 and silently drop the trailing zero / formatting. Marking
# the cell as Text first keeps the assignment literal.
$textForceCells = @(
    "D5",
    "D6",
    "D7",
    "D9",
    "D10",
    "D11",
    "D12",
    "D14",
    "D17",
    "D19",
    "D20",
    "D21",
    "D23",
    "D24",
    "D27",
    "D29",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D46",
    "D47",
    "D49",
    "D50",
    "D51",
)
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "48.173.13"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").Value = "2.494.53"
$ws.Range("E3").Value = "  -1.43%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "316.96"
$ws.Range("E5").Value = "  -2.02%  "
$ws.Range("D6").Value = "105.71"
$ws.Range("E6").Value = "  -3.29%  "
$ws.Range("D7").Value = "0.518"
$ws.Range("E7").Value = "  -1.90%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "0.534"
$ws.Range("E9").Value = "  -4.40%  "
$ws.Range("D10").Value = "38.77"
$ws.Range("E10").Value = "  -4.61%  "
$ws.Range("D11").Value = "20.28"
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("D12").Value = "0.0802"
$ws.Range("E12").Value = "  -2.46%  "
$ws.Range("E13").Value = "  -0.45%  "
$ws.Range("D14").Value = "7.06"
$ws.Range("E14").Value = "  -3.12%  "
$ws.Range("D15").Value = "2.881.58"
$ws.Range("E15").Value = "  -1.49%  "
$ws.Range("D16").Value = "2.494.23"
$ws.Range("D17").Value = "0.828"
$ws.Range("E17").Value = "  -3.83%  "
$ws.Range("D18").Value = "47.968.27"
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("D19").Value = "2.97"
$ws.Range("E19").Value = "  +8.92%  "
$ws.Range("D20").Value = "12.85"
$ws.Range("E20").Value = "  -3.60%  "
$ws.Range("D21").Value = "6.58"
$ws.Range("E21").Value = "  -0.82%  "
$ws.Range("D22").Value = "0.0₃0930"
$ws.Range("E22").Value = "  -1.75%  "
$ws.Range("D23").Value = "70.98"
$ws.Range("E23").Value = "  -2.07%  "
$ws.Range("D24").Value = "268.80"
$ws.Range("E24").Value = "  -0.61%  "
$ws.Range("E25").Value = "  -3.30%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "25.77"
$ws.Range("E27").Value = "  -2.03%  "
$ws.Range("E28").Value = "  +4.42%  "
$ws.Range("D29").Value = "9.70"
$ws.Range("E29").Value = "  -4.75%  "
$ws.Range("E30").Value = "  -5.44%  "
$ws.Range("D31").Value = "34.47"
$ws.Range("E31").Value = "  -2.75%  "
$ws.Range("D32").Value = "49.34"
$ws.Range("E32").Value = "  -0.87%  "
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "5.29"
$ws.Range("E34").Value = "  -2.28%  "
$ws.Range("B35").Value = "Celestia"
$ws.Range("C35").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D35").Value = "19.00"
$ws.Range("E35").Value = "  -5.14%  "
$ws.Range("D36").Value = "0.0772"
$ws.Range("E36").Value = "  -2.46%  "
$ws.Range("D37").Value = "1.93"
$ws.Range("E37").Value = "  -3.13%  "
$ws.Range("D38").Value = "4.57"
$ws.Range("E38").Value = "  -3.78%  "
$ws.Range("D39").Value = "2.87"
$ws.Range("E39").Value = "  -4.60%  "
$ws.Range("D40").Value = "122.23"
$ws.Range("E40").Value = "  +3.26%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").Value = "22.43"
$ws.Range("E41").Value = "  -0.80%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "0.110"
$ws.Range("E42").Value = "  -1.77%  "
$ws.Range("D43").Value = "2.22"
$ws.Range("E43").Value = "  +1.08%  "
$ws.Range("D44").Value = "0.0301"
$ws.Range("E44").Value = "  +0.65%  "
$ws.Range("D45").Value = "2.000.24"
$ws.Range("E45").Value = "  -0.59%  "
$ws.Range("D46").Value = "3.12"
$ws.Range("E46").Value = "  -1.24%  "
$ws.Range("D47").Value = "1.89"
$ws.Range("E47").Value = "  -0.42%  "
$ws.Range("E48").Value = "  -1.35%  "
$ws.Range("D49").Value = "8.90"
$ws.Range("E49").Value = "  -2.56%  "
$ws.Range("D50").Value = "5.15"
$ws.Range("E50").Value = "  -2.07%  "
$ws.Range("D51").Value = "78.59"
$ws.Range("E51").Value = "  -2.09%  "

# Put the cells back on the workbook default style now that the literal text
# is committed, so no stray NumberFormat/style survives the edit.
foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
